# Week 13 results update
# - Removes the "Image" column values from the old "Week 12" rows (134-145)
# - Appends a new "Week 13" block (rows 146-157) with the same 12 teams,
#   new Points For / Points Against values, and the Image column values
#   that used to sit on the Week 12 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Clear the Image column on the old "Week 12" rows (134:145) -- the
#    Image column now belongs to the newest week only.
$ws.Range("E134:E145").ClearContents()

# 2) Append the new Week 13 rows (146-157), carrying the per-team icon
#    path down from where it used to live on the Week 12 rows. Inserting
#    each row (rather than just writing past the used range) makes Excel
#    carry the Points For / Points Against number formatting down from
#    the row above, same as the rest of the table.
for ($r = 146; $r -le 157; $r++) {
  $ws.Rows.Item($r).Insert(-4121)
}

$week13 = "Week 13"

$ws.Cells.Item(146,1).Value = "Kauaireek Hill"
$ws.Cells.Item(146,2).Value = $week13
$ws.Cells.Item(146,3).Value = 119.14
$ws.Cells.Item(146,4).Value = 81.62
$ws.Cells.Item(146,5).Value = "Team Icons/kauaireek-modified.png"

$ws.Cells.Item(147,1).Value = "Chasing dank Herb"
$ws.Cells.Item(147,2).Value = $week13
$ws.Cells.Item(147,3).Value = 102.1
$ws.Cells.Item(147,4).Value = 98.72
$ws.Cells.Item(147,5).Value = "Team Icons/chasing-modified.png"

$ws.Cells.Item(148,1).Value = "Ju Ju Smith Poopster"
$ws.Cells.Item(148,2).Value = $week13
$ws.Cells.Item(148,3).Value = 89.74
$ws.Cells.Item(148,4).Value = 74.5
$ws.Cells.Item(148,5).Value = "Team Icons/juju-modified.png"

$ws.Cells.Item(149,1).Value = "Bye Breece See You in ValHalla"
$ws.Cells.Item(149,2).Value = $week13
$ws.Cells.Item(149,3).Value = 115.3
$ws.Cells.Item(149,4).Value = 140.96
$ws.Cells.Item(149,5).Value = "Team Icons/breece-modified.png"

$ws.Cells.Item(150,1).Value = "Cooking with Gas"
$ws.Cells.Item(150,2).Value = $week13
$ws.Cells.Item(150,3).Value = 140.96
$ws.Cells.Item(150,4).Value = 115.3
$ws.Cells.Item(150,5).Value = "Team Icons/cooking-modified.png"

$ws.Cells.Item(151,1).Value = "Dulcich de Leche"
$ws.Cells.Item(151,2).Value = $week13
$ws.Cells.Item(151,3).Value = 106.78
$ws.Cells.Item(151,4).Value = 120.6
$ws.Cells.Item(151,5).Value = "Team Icons/dulcich-modified.png"

$ws.Cells.Item(152,1).Value = "Dillon Panthers"
$ws.Cells.Item(152,2).Value = $week13
$ws.Cells.Item(152,3).Value = 123.72
$ws.Cells.Item(152,4).Value = 111.86
$ws.Cells.Item(152,5).Value = "Team Icons/dillon-modified.png"

$ws.Cells.Item(153,1).Value = "Daemon and the Rightful Heirs"
$ws.Cells.Item(153,2).Value = $week13
$ws.Cells.Item(153,3).Value = 111.86
$ws.Cells.Item(153,4).Value = 123.72
$ws.Cells.Item(153,5).Value = "Team Icons/daemon-modified.png"

$ws.Cells.Item(154,1).Value = "Krombopulos Michael Evans"
$ws.Cells.Item(154,2).Value = $week13
$ws.Cells.Item(154,3).Value = 74.5
$ws.Cells.Item(154,4).Value = 89.74
$ws.Cells.Item(154,5).Value = "Team Icons/krombopulos-modified.png"

$ws.Cells.Item(155,1).Value = "Freier Freier Pants on Fire"
$ws.Cells.Item(155,2).Value = $week13
$ws.Cells.Item(155,3).Value = 81.62
$ws.Cells.Item(155,4).Value = 119.14
$ws.Cells.Item(155,5).Value = "Team Icons/freier-modified.png"

$ws.Cells.Item(156,1).Value = "Christian Kirk Cousins"
$ws.Cells.Item(156,2).Value = $week13
$ws.Cells.Item(156,3).Value = 98.72
$ws.Cells.Item(156,4).Value = 102.1
$ws.Cells.Item(156,5).Value = "Team Icons/sir-modified.png"

$ws.Cells.Item(157,1).Value = "Fantasy Football Champion 2022"
$ws.Cells.Item(157,2).Value = $week13
$ws.Cells.Item(157,3).Value = 120.6
$ws.Cells.Item(157,4).Value = 106.78
$ws.Cells.Item(157,5).Value = "Team Icons/fantasy-modified.png"

# 3) Leave the cursor where a user would land after typing the last row.
[void]$ws.Range("C158").Select()
